# Edit script: turns the "Dark Matter" essay into the "Government" essay,
# renames the author and e-mail, rewrites the body/summary text, and adds a
# trailing empty paragraph - matching the supplied OOXML diff.

$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2
$br = [char]11   # manual line break (renders as <w:br/>)

# ---------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Beyond the Pale Horizon: Exploring the Mysteries of Dark Matter", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Navigating the Complexities of Government: A Primer for High School Students", $wdReplaceAll) | Out-Null

# ---------------------------------------------------------------------
# 2. Author name
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Sarah Walker", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Olivia Harrison", $wdReplaceAll) | Out-Null

# ---------------------------------------------------------------------
# 3. E-mail line: "walker" + "." + "sarah@astroscience" + "." + "org"
#    becomes "oliviaharrison@valid" + "." + "com"
#    (the middle "." run is kept as-is; the other two text runs are dropped)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "sarah@astroscience.", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "", $wdReplaceAll) | Out-Null

$d.Content.Find.Execute(
    "walker", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "oliviaharrison@valid", $wdReplaceAll) | Out-Null

$d.Content.Find.Execute(
    "org", $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "com", $wdReplaceAll) | Out-Null

# ---------------------------------------------------------------------
# 4. Intro / body paragraph (the long paragraph right after the blank line)
# ---------------------------------------------------------------------
$introPara = $d.Paragraphs.Item(5)
$introRange = $d.Range($introPara.Range.Start, $introPara.Range.End - 1)

$introText = (
    "Governments, with their intricate systems and dynamic processes, play a pivotal role in shaping our societies." +
    " From local municipalities to global alliances, governance structures influence our daily lives in myriad ways." +
    " This essay delves into the fascinating realm of government, providing high school students with a comprehensive understanding of its diverse forms, fundamental principles, and real-world applications." +
    " We will venture into the historical evolution of governance, exploring the theories, ideologies, and institutions that have shaped the modern political landscape." +
    $br + $br +
    "In the tapestry of human history, governance has been a constant companion, evolving through various forms and structures." +
    " From ancient civilizations to modern nation-states, the concept of government has undergone a remarkable transformation." +
    " We will trace the journey of governance, from its early roots in tribal societies to the complex systems we see today." +
    " By examining historical case studies, we will uncover the factors that have influenced the development of different forms of government, exploring the interplay between power, authority, and legitimacy." +
    $br + $br +
    "Governments exist in a kaleidoscope of forms, each with its unique characteristics and implications." +
    " As we dive into the study of government, we will encounter a range of political systems, including democracies, autocracies, monarchies, and more." +
    " We will analyze the structures, processes, and institutions that define these systems, comparing their strengths and weaknesses." +
    " By examining the different branches of government, such as the executive, legislative, and judicial, we will gain insights into the intricate web of checks and balances that ensure accountability and prevent the concentration of power." +
    $br + $br +
    "Body:" +
    $br +
    "The principles that underpin governments are as multifaceted as the systems themselves." +
    " Justice, equality, freedom, and security are among the core values that governments strive to uphold." +
    " We will delve into the philosophical and historical foundations of these principles, exploring how they have shaped the evolution of governance." +
    " By examining landmark legal cases, constitutional frameworks, and international treaties, we will uncover the intricate relationship between government actions and the rights and responsibilities of citizens." +
    $br + $br +
    "The practical applications of government are vast and encompass a wide range of areas that directly impact our daily lives." +
    " From the provision of essential services like healthcare and education to the regulation of businesses and industries, governments play a crucial role in shaping our societies." +
    " We will investigate the various functions of government, exploring how policies are formulated, implemented, and evaluated." +
    " By examining real-world examples, we will gain a deeper understanding of the complex challenges and opportunities that governments face in addressing the needs of their citizens." +
    $br + $br +
    "The study of government is incomplete without acknowledging the role of citizens in the governance process." +
    " As active participants in a democracy, citizens have a responsibility to engage with their government, hold it accountable, and contribute to decision-making processes." +
    " We will explore the various avenues through which citizens can exercise their rights and influence government actions." +
    " By examining case studies of successful citizen engagement, we will uncover the transformative power of active participation in shaping a more just and equitable society."
)

$introRange.Text = $introText

# ---------------------------------------------------------------------
# 5. Summary paragraph (last paragraph of the document)
# ---------------------------------------------------------------------
$summaryPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$summaryRange = $d.Range($summaryPara.Range.Start, $summaryPara.Range.End - 1)

$summaryText = (
    "In this essay, we embarked on a journey through the captivating world of government." +
    " We explored the historical evolution of governance, examining the diverse forms, structures, and principles that define political systems across the globe." +
    " We delved into the practical applications of government, investigating how policies are formulated and implemented, and how they impact the lives of citizens." +
    " Finally, we recognized the crucial role of citizens in the governance process, emphasizing their rights, responsibilities, and the transformative power of active participation." +
    " Throughout this exploration, we gained a deeper appreciation for the complexity and significance of government in shaping our societies."
)

$summaryRange.Text = $summaryText

# ---------------------------------------------------------------------
# 6. Append a trailing empty paragraph at the very end of the document
# ---------------------------------------------------------------------
$endPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endPara.Range.InsertParagraphAfter() | Out-Null
